$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.536.08'
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").Value = '1.472.06'
$ws.Range("E3").Value = '  +1.90%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("E5").Value = '  +4.98%  '

$ws.Range("D6").Value = '''277.67'
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").Value = '''0.3616'
$ws.Range("E7").Value = '  -1.34%  '

$ws.Range("D8").Value = '''0.3066'
$ws.Range("E8").Value = '  -2.18%  '

$ws.Range("E9").Value = '  +1.53%  '

$ws.Range("D10").Value = '''1.069'
$ws.Range("E10").Value = '  +4.58%  '

$ws.Range("D11").Value = '''0.06634'
$ws.Range("E11").Value = '  +1.43%  '

$ws.Range("E12").Value = '  +0.11%  '

$ws.Range("D13").Value = '''5.522'
$ws.Range("E13").Value = '  +2.21%  '

$ws.Range("D14").Value = '''18.09'
$ws.Range("E14").Value = '  +2.73%  '

$ws.Range("B15").Value = 'Dai'
$ws.Range("C15").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D15").Value = '''0.9578'
$ws.Range("E15").Value = '  +2.38%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''6.159'
$ws.Range("E16").Value = '  +1.36%  '

$ws.Range("E17").Value = '  +0.90%  '

$ws.Range("D18").Value = '1.473.72'
$ws.Range("E18").Value = '  +2.32%  '

$ws.Range("D19").Value = '''0.05921'
$ws.Range("E19").Value = '  +5.11%  '

$ws.Range("D20").Value = '''68.92'
$ws.Range("E20").Value = '  +1.87%  '

$ws.Range("D21").Value = '''5.490'
$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").Value = '''11.16'
$ws.Range("E23").Value = '  +3.05%  '

$ws.Range("D24").Value = '''2.260'
$ws.Range("E24").Value = '  +0.78%  '

$ws.Range("D25").Value = '20.542.74'
$ws.Range("E25").Value = '  +1.59%  '

$ws.Range("E26").Value = '  +4.94%  '

$ws.Range("E27").Value = '  -3.13%  '

$ws.Range("D28").Value = '''17.13'
$ws.Range("E28").Value = '  +0.94%  '

$ws.Range("D29").Value = '1.636.29'
$ws.Range("E29").Value = '  +2.69%  '

$ws.Range("D30").Value = '''113.49'
$ws.Range("E30").Value = '  +2.84%  '

$ws.Range("D31").Value = '''3.916'
$ws.Range("E31").Value = '  +3.38%  '

$ws.Range("D32").Value = '''4.963'
$ws.Range("E32").Value = '  +2.41%  '

$ws.Range("D33").Value = '''0.07995'

$ws.Range("D34").Value = '''0.8053'
$ws.Range("E34").Value = '  -0.54%  '

$ws.Range("D35").Value = '''1.514'
$ws.Range("E35").Value = '  +4.63%  '

$ws.Range("D36").Value = '''1.214'
$ws.Range("E36").Value = '  +4.97%  '

$ws.Range("D37").Value = '''0.05742'
$ws.Range("E37").Value = '  -4.49%  '

$ws.Range("D38").Value = '''4.720'
$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("E39").Value = '  +3.06%  '

$ws.Range("D40").Value = '''0.9583'
$ws.Range("E40").Value = '  +2.82%  '

$ws.Range("D41").Value = '''10.35'
$ws.Range("E41").Value = '  +1.38%  '

$ws.Range("D42").Value = '''0.1874'
$ws.Range("E42").Value = '  +2.04%  '

$ws.Range("E43").Value = '  +3.98%  '

$ws.Range("D44").Value = '''0.5270'
$ws.Range("E44").Value = '  +0.52%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '''3.524'
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''12.13'
$ws.Range("E46").Value = '  +0.44%  '

$ws.Range("D47").Value = '''117.95'
$ws.Range("E47").Value = '  -0.77%  '

$ws.Range("D48").Value = '''0.5194'
$ws.Range("E48").Value = '  +0.78%  '

$ws.Range("D49").Value = '''1.811'
$ws.Range("E49").Value = '  +2.37%  '

$ws.Range("D50").Value = '''0.06465'
$ws.Range("E50").Value = '  +2.15%  '

$ws.Range("D51").Value = '''0.9860'
$ws.Range("E51").Value = '  -0.65%  '

Write-Host "Cryptos list updated"